$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as text so that values such as
# "446.65" or "0.0000320" are not auto-converted to floating point numbers by
# the numeric auto-detection in Range.Value. We flip the format to Text,
# assign all values, then restore the original (Normal) style so the saved
# file keeps the same cell styling as before the edit.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "66.819.13"
$ws.Range("E2").Value = "  +1.10%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "3.808.50"
$ws.Range("E3").Value = "  +0.83%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 (BNB)
$ws.Range("D5").Value = "446.65"
$ws.Range("E5").Value = "  +6.71%  "

# Row 6 (Solana)
$ws.Range("D6").Value = "145.74"
$ws.Range("E6").Value = "  +14.56%  "

# Row 7 (XRP)
$ws.Range("D7").Value = "0.623"
$ws.Range("E7").Value = "  +4.69%  "

# Row 8 (USDC)
$ws.Range("E8").Value = "  -0.06%  "

# Row 9 (Cardano)
$ws.Range("D9").Value = "0.734"
$ws.Range("E9").Value = "  +2.98%  "

# Row 10 (Dogecoin)
$ws.Range("D10").Value = "0.155"
$ws.Range("E10").Value = "  -2.63%  "

# Row 11 (ShibaInu)
$ws.Range("D11").Value = "0.0000320"
$ws.Range("E11").Value = "  -6.15%  "

# Row 12 (Avalanche)
$ws.Range("D12").Value = "43.36"
$ws.Range("E12").Value = "  +9.78%  "

# Row 13 (Polkadot)
$ws.Range("D13").Value = "10.31"
$ws.Range("E13").Value = "  +3.10%  "

# Row 14 (WrappedliquidstakedEther2.0)
$ws.Range("D14").Value = "4.412.93"
$ws.Range("E14").Value = "  +0.97%  "

# Row 15 (Uniswap)
$ws.Range("D15").Value = "15.01"
$ws.Range("E15").Value = "  -6.08%  "

# Row 16 (WrappedEther)
$ws.Range("D16").Value = "3.865.91"
$ws.Range("E16").Value = "  +2.22%  "

# Row 17 (TRON)
$ws.Range("E17").Value = "  -0.19%  "

# Row 18 (Chainlink)
$ws.Range("D18").Value = "19.91"
$ws.Range("E18").Value = "  +3.18%  "

# Row 19
$ws.Range("E19").Value = "  +6.79%  "

# Row 20
$ws.Range("D20").Value = "66.912.71"
$ws.Range("E20").Value = "  +1.15%  "

# Row 21
$ws.Range("D21").Value = "422.70"
$ws.Range("E21").Value = "  +5.01%  "

# Row 22
$ws.Range("D22").Value = "14.63"
$ws.Range("E22").Value = "  +3.78%  "

# Row 23
$ws.Range("D23").Value = "3.23"
$ws.Range("E23").Value = "  +9.19%  "

# Row 24
$ws.Range("D24").Value = "86.62"
$ws.Range("E24").Value = "  +4.62%  "

# Row 25
$ws.Range("D25").Value = "37.24"
$ws.Range("E25").Value = "  +1.92%  "

# Row 26
$ws.Range("D26").Value = "3.43"
$ws.Range("E26").Value = "  +8.12%  "

# Row 27
$ws.Range("D27").Value = "5.51"
$ws.Range("E27").Value = "  -3.31%  "

# Row 28 (Filecoin -> RenderToken)
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "9.47"
$ws.Range("E28").Value = "  +19.24%  "

# Row 29 (RenderToken -> Filecoin)
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "9.66"
$ws.Range("E29").Value = "  +4.59%  "

# Row 30
$ws.Range("D30").Value = "747.65"
$ws.Range("E30").Value = "  +7.43%  "

# Row 31
$ws.Range("D31").Value = "13.66"
$ws.Range("E31").Value = "  +12.43%  "

# Row 32
$ws.Range("E32").Value = "  +11.89%  "

# Row 33
$ws.Range("E33").Value = "  -1.08%  "

# Row 34
$ws.Range("D34").Value = "43.05"
$ws.Range("E34").Value = "  +15.91%  "

# Row 35
$ws.Range("D35").Value = "0.155"
$ws.Range("E35").Value = "  +4.24%  "

# Row 36
$ws.Range("D36").Value = "58.56"
$ws.Range("E36").Value = "  +7.71%  "

# Row 37
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.18%  "

# Row 38
$ws.Range("D38").Value = "5.47"
$ws.Range("E38").Value = "  +17.55%  "

# Row 39
$ws.Range("D39").Value = "0.0474"
$ws.Range("E39").Value = "  +5.81%  "

# Row 40
$ws.Range("D40").Value = "0.349"
$ws.Range("E40").Value = "  +19.86%  "

# Row 41
$ws.Range("D41").Value = "2.87"
$ws.Range("E41").Value = "  -1.39%  "

# Row 42
$ws.Range("E42").Value = "  +0.08%  "

# Row 43 (Stellar -> PEPE)
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").Value = "0.0₃0672"
$ws.Range("E43").Value = "  -10.77%  "

# Row 44 (PEPE -> Stellar)
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "0.140"
$ws.Range("E44").Value = "  +4.92%  "

# Row 45
$ws.Range("E45").Value = "  +5.64%  "

# Row 46
$ws.Range("E46").Value = "  +3.74%  "

# Row 47
$ws.Range("E47").Value = "  +11.66%  "

# Row 48
$ws.Range("E48").Value = "  +2.23%  "

# Row 49
$ws.Range("D49").Value = "2.10"
$ws.Range("E49").Value = "  +4.58%  "

# Row 50
$ws.Range("D50").Value = "2.65"
$ws.Range("E50").Value = "  +6.62%  "

# Row 51
$ws.Range("D51").Value = "2.86"
$ws.Range("E51").Value = "  +5.21%  "

# Restore the original (unstyled / Normal) style for the whole range so that
# no stray formatting differences are introduced relative to the source file.
$priceVolRange.Style = "Normal"
